$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from P1 into the new header cells Q1:T1
$ws.Range("P1").Copy() | Out-Null
$ws.Range("Q1:T1").PasteSpecial(-4122) | Out-Null

# Header row (row 1): new "station 6 / station 7" columns
$ws.Range("Q1").Value = "Estación más cercana 6"
$ws.Range("R1").Value = "Estación más cercana 7"
$ws.Range("S1").Value = "Inicio estación más cercana 6"
$ws.Range("T1").Value = "Inicio estación más cercana 7"

# Data rows 2-46
$ws.Range("Q2").Value = "PATCX"
$ws.Range("R2").Value = "PX02"
$ws.Range("S2").Value = "2014-03-16T21:16:41"
$ws.Range("T2").Value = "2014-03-16T21:16:42"
$ws.Range("Q3").Value = "PB11"
$ws.Range("R3").Value = "PX03"
$ws.Range("S3").Value = "2014-03-17T05:11:48"
$ws.Range("T3").Value = "2014-03-17T05:11:49"
$ws.Range("Q4").Value = "PB11"
$ws.Range("R4").Value = "PATCX"
$ws.Range("S4").Value = "2014-03-22T13:00:13"
$ws.Range("T4").Value = "2014-03-22T13:00:14"
$ws.Range("Q5").Value = "PB11"
$ws.Range("R5").Value = "PB12"
$ws.Range("S5").Value = "2014-03-23T18:20:15"
$ws.Range("T5").Value = "2014-03-23T18:20:17"
$ws.Range("Q6").Value = "TA01"
$ws.Range("R6").Value = "PB12"
$ws.Range("S6").Value = "2014-04-01T23:47:02"
$ws.Range("T6").Value = "2014-04-01T23:47:02"
$ws.Range("Q7").Value = "PX03"
$ws.Range("R7").Value = "PB11"
$ws.Range("S7").Value = "2014-04-03T01:58:40"
$ws.Range("T7").Value = "2014-04-03T01:58:42"
$ws.Range("Q8").Value = "PB02"
$ws.Range("R8").Value = "PB01"
$ws.Range("S8").Value = "2014-04-03T02:43:26"
$ws.Range("T8").Value = "2014-04-03T02:43:28"
$ws.Range("Q9").Value = "PX03"
$ws.Range("R9").Value = "PB01"
$ws.Range("S9").Value = "2014-04-03T05:26:25"
$ws.Range("T9").Value = "2014-04-03T05:26:27"
$ws.Range("Q10").Value = "PX03"
$ws.Range("R10").Value = "PB01"
$ws.Range("S10").Value = "2014-04-04T01:38:03"
$ws.Range("T10").Value = "2014-04-04T01:38:05"
$ws.Range("Q11").Value = "PB11"
$ws.Range("R11").Value = "PX02"
$ws.Range("S11").Value = "2014-04-01T23:58:13"
$ws.Range("T11").Value = "2014-04-01T23:58:13"
$ws.Range("Q12").Value = "PX03"
$ws.Range("R12").Value = "PB01"
$ws.Range("S12").Value = "2014-04-11T00:01:56"
$ws.Range("T12").Value = "2014-04-11T00:01:58"
$ws.Range("Q13").Value = "MT02"
$ws.Range("R13").Value = "CO04"
$ws.Range("S13").Value = "2014-08-23T22:32:31"
$ws.Range("T13").Value = "2014-08-23T22:32:33"
$ws.Range("Q14").Value = "VA01"
$ws.Range("R14").Value = "ROC1"
$ws.Range("S14").Value = "2015-09-16T22:54:51"
$ws.Range("T14").Value = "2015-09-16T22:54:52"
$ws.Range("Q15").Value = "VA03"
$ws.Range("R15").Value = "GO04"
$ws.Range("S15").Value = "2015-09-16T23:19:00"
$ws.Range("T15").Value = "2015-09-16T23:19:01"
$ws.Range("Q16").Value = "VA01"
$ws.Range("R16").Value = "ROC1"
$ws.Range("S16").Value = "2015-09-16T23:16:28"
$ws.Range("T16").Value = "2015-09-16T23:16:30"
$ws.Range("Q17").Value = "GO04"
$ws.Range("R17").Value = "TLL"
$ws.Range("S17").Value = "2015-09-17T01:41:21"
$ws.Range("T17").Value = "2015-09-17T01:41:21"
$ws.Range("Q18").Value = "GO04"
$ws.Range("R18").Value = "TLL"
$ws.Range("S18").Value = "2015-09-17T03:55:35"
$ws.Range("T18").Value = "2015-09-17T03:55:35"
$ws.Range("Q19").Value = "VA01"
$ws.Range("R19").Value = "GO04"
$ws.Range("S19").Value = "2015-09-17T04:10:47"
$ws.Range("T19").Value = "2015-09-17T04:10:48"
$ws.Range("Q20").Value = "MT07"
$ws.Range("R20").Value = "MT02"
$ws.Range("S20").Value = "2015-09-18T09:11:00"
$ws.Range("T20").Value = "2015-09-18T09:11:01"
$ws.Range("Q21").Value = "GO04"
$ws.Range("R21").Value = "TLL"
$ws.Range("S21").Value = "2015-09-19T05:07:03"
$ws.Range("T21").Value = "2015-09-19T05:07:03"
$ws.Range("Q22").Value = "MT07"
$ws.Range("R22").Value = "MT02"
$ws.Range("S22").Value = "2015-09-19T12:52:35"
$ws.Range("T22").Value = "2015-09-19T12:52:36"
$ws.Range("Q23").Value = "VA01"
$ws.Range("R23").Value = "ROC1"
$ws.Range("S23").Value = "2015-09-21T05:39:53"
$ws.Range("T23").Value = "2015-09-21T05:39:55"
$ws.Range("Q24").Value = "ROC1"
$ws.Range("R24").Value = "MT07"
$ws.Range("S24").Value = "2015-09-21T17:40:17"
$ws.Range("T24").Value = "2015-09-21T17:40:17"
$ws.Range("Q25").Value = "GO04"
$ws.Range("R25").Value = "TLL"
$ws.Range("S25").Value = "2015-09-22T07:13:18"
$ws.Range("T25").Value = "2015-09-22T07:13:18"
$ws.Range("Q26").Value = "IN40"
$ws.Range("R26").Value = "IN41"
$ws.Range("S26").Value = "2015-09-26T02:51:29"
$ws.Range("T26").Value = "2015-09-26T02:51:29"
$ws.Range("Q27").Value = "GO04"
$ws.Range("R27").Value = "TLL"
$ws.Range("S27").Value = "2015-11-07T07:04:51"
$ws.Range("T27").Value = "2015-11-07T07:04:51"
$ws.Range("Q28").Value = "GO04"
$ws.Range("R28").Value = "TLL"
$ws.Range("S28").Value = "2015-11-11T01:54:55"
$ws.Range("T28").Value = "2015-11-11T01:54:55"
$ws.Range("Q29").Value = "GO04"
$ws.Range("R29").Value = "TLL"
$ws.Range("S29").Value = "2015-11-11T02:46:36"
$ws.Range("T29").Value = "2015-11-11T02:46:36"
$ws.Range("Q30").Value = "PB19"
$ws.Range("R30").Value = "PB15"
$ws.Range("S30").Value = "2015-11-27T21:00:42"
$ws.Range("T30").Value = "2015-11-27T21:00:48"
$ws.Range("Q31").Value = "TLL"
$ws.Range("R31").Value = "CO02"
$ws.Range("S31").Value = "2016-02-10T00:33:15"
$ws.Range("T31").Value = "2016-02-10T00:33:16"
$ws.Range("Q32").Value = "VA01"
$ws.Range("R32").Value = "IN41"
$ws.Range("S32").Value = "2015-09-17T04:10:47"
$ws.Range("T32").Value = "2015-09-17T04:10:48"
$ws.Range("Q33").Value = "GO04"
$ws.Range("R33").Value = "TLL"
$ws.Range("S33").Value = "2015-09-22T07:13:18"
$ws.Range("T33").Value = "2015-09-22T07:13:18"
$ws.Range("Q34").Value = "VA01"
$ws.Range("R34").Value = "ROC1"
$ws.Range("S34").Value = "2015-09-21T05:39:53"
$ws.Range("T34").Value = "2015-09-21T05:39:54"
$ws.Range("Q35").Value = "MT07"
$ws.Range("R35").Value = "MT02"
$ws.Range("S35").Value = "2015-09-19T12:52:35"
$ws.Range("T35").Value = "2015-09-19T12:52:36"
$ws.Range("Q36").Value = "ROC1"
$ws.Range("R36").Value = "MT07"
$ws.Range("S36").Value = "2017-04-23T02:36:18"
$ws.Range("T36").Value = "2017-04-23T02:36:18"
$ws.Range("Q37").Value = "ROC1"
$ws.Range("R37").Value = "MT07"
$ws.Range("S37").Value = "2017-04-24T21:38:41"
$ws.Range("T37").Value = "2017-04-24T21:38:41"
$ws.Range("Q38").Value = "CO06"
$ws.Range("R38").Value = "CO10"
$ws.Range("S38").Value = "2019-01-20T01:32:59"
$ws.Range("T38").Value = "2019-01-20T01:33:02"
$ws.Range("Q39").Value = "GO04"
$ws.Range("R39").Value = "TLL"
$ws.Range("S39").Value = "2019-06-14T00:19:27"
$ws.Range("T39").Value = "2019-06-14T00:19:27"
$ws.Range("Q40").Value = "ROC1"
$ws.Range("R40").Value = "CO06"
$ws.Range("S40").Value = "2019-11-04T21:53:41"
$ws.Range("T40").Value = "2019-11-04T21:53:41"
$ws.Range("Q41").Value = "PB16"
$ws.Range("R41").Value = "PSGCX"
$ws.Range("S41").Value = "2019-12-03T08:46:51"
$ws.Range("T41").Value = "2019-12-03T08:46:51"
$ws.Range("Q42").Value = "AC05"
$ws.Range("R42").Value = "CO10"
$ws.Range("S42").Value = "2020-09-01T04:09:45"
$ws.Range("T42").Value = "2020-09-01T04:09:45"
$ws.Range("Q43").Value = "AC05"
$ws.Range("R43").Value = "CO10"
$ws.Range("S43").Value = "2020-09-01T04:30:18"
$ws.Range("T43").Value = "2020-09-01T04:30:18"
$ws.Range("Q44").Value = "AC05"
$ws.Range("R44").Value = "CO10"
$ws.Range("S44").Value = "2020-09-01T21:09:35"
$ws.Range("T44").Value = "2020-09-01T21:09:35"
$ws.Range("Q45").Value = "TLL"
$ws.Range("R45").Value = "CO03"
$ws.Range("S45").Value = "2020-09-06T01:17:06"
$ws.Range("T45").Value = "2020-09-06T01:17:09"
$ws.Range("Q46").Value = "PB03"
$ws.Range("R46").Value = "PB09"
$ws.Range("S46").Value = "2020-09-11T07:36:06"
$ws.Range("T46").Value = "2020-09-11T07:36:07"
